# Kanban board update: "Account System: Change Password" task moves from the
# "Not Started" column (A) to the "Done" column (C) — i.e. the task was
# completed, so it is removed from its spot in column A (rows below it in
# column A shift up by one) and appended to the first empty cell in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$taskName = "Account System: Change Password"

# Shift column A items up starting at row 5 (overwrite A5 with A6, A6 with A7, ...)
for ($r = 5; $r -le 9; $r++) {
    $nextVal = $ws.Cells.Item($r + 1, 1).Value2
    $ws.Cells.Item($r, 1).Value = $nextVal
}

# The last row that used to hold a value (row 10) is now empty.
$ws.Cells.Item(10, 1).ClearContents()

# Append the completed task to the first empty row in the "Done" column (C15).
$ws.Range("C15").Value = $taskName
